$wb = $excel.ActiveWorkbook

# --- Update data on "GHEDT (HYTS)" sheet ---
$ws2 = $wb.Worksheets.Item("GHEDT (HYTS)")

$ws2.Range("B2").Value = 128.15841855
$ws2.Range("C2").Value = 131.42350042
$ws2.Range("D2").Value = 131.11602901

$ws2.Range("B3").Value = 136.6598564
$ws2.Range("C3").Value = 130.12624536
$ws2.Range("D3").Value = 129.79844448

$ws2.Range("B4").Value = 139.43443151
$ws2.Range("C4").Value = 131.62053447
$ws2.Range("D4").Value = 127.49687689

$ws2.Range("B5").Value = 141.38486327
$ws2.Range("C5").Value = 132.71390552
$ws2.Range("D5").Value = 125.78451994

$ws2.Range("B6").Value = 142.91893499
$ws2.Range("C6").Value = 133.52339467
$ws2.Range("D6").Value = 125.48150927

# --- Update data on "GHEDT (HOTS)" sheet ---
$ws3 = $wb.Worksheets.Item("GHEDT (HOTS)")

$ws3.Range("B2").Value = 129.201
$ws3.Range("C2").Value = 137.8979
$ws3.Range("D2").Value = 148.915

$ws3.Range("B3").Value = 136.8608
$ws3.Range("C3").Value = 131.0179
$ws3.Range("D3").Value = 133.1613

$ws3.Range("B4").Value = 139.6402
$ws3.Range("C4").Value = 132.2753
$ws3.Range("D4").Value = 127.6939

$ws3.Range("B5").Value = 141.5259
$ws3.Range("C5").Value = 133.2482
$ws3.Range("D5").Value = 126.116

$ws3.Range("B6").Value = 142.9455
$ws3.Range("C6").Value = 134.1695
$ws3.Range("D6").Value = 125.7789

# Row 6 on "GHEDT (HYTS)" grows from a 12.8pt row to the sheet's standard 13.8pt row
$ws2.Rows.Item(6).RowHeight = 13.8

# --- Update selection / active cell on each sheet ---
# GLHEPRO: no longer the selected tab
$ws1 = $wb.Worksheets.Item("GLHEPRO")
$ws1.Range("A1").Select()

# GHEDT (HYTS): becomes the active/selected tab, active cell A1
$ws2.Select()
$ws2.Range("A1").Select()

# GHEDT (HOTS): active cell moves to A1
$ws3.Range("A1").Select()

# Make GHEDT (HYTS) the active sheet/tab for the workbook
$ws2.Activate()
